$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: fff / fff / evolutic
$ws.Range("A2").Value = "fff"

# Row 4: aaa / evolutic (B4, C4 first)
$ws.Range("B4").Value = "aaa"
$ws.Range("C4").Value = "evolutic"

# Row 5: EEE / ccc
$ws.Range("A5").Value = "EEE"
$ws.Range("B5").Value = "ccc"

# Row 6: FFF / xxx
$ws.Range("A6").Value = "FFF"
$ws.Range("B6").Value = "xxx"

# Row 3: jaqussd / Jaquqssde / evolutic
$ws.Range("A3").Value = "jaqussd"
$ws.Range("C3").Value = "evolutic"

# Row 4: vvvvv (A4 filled last)
$ws.Range("A4").Value = "vvvvv"

$ws.Range("A5").Select()
